$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "www.ulutasmedicaljournal.com",
    "experimentalbiomedicalresearch.com",
    "medicaljournal.gazi.edu.tr",
    "beslenmevediyetdergisi.org",
    "www.jsoah.com",
    "actamedica.org",
    "www.cityhealthj.org",
    "jointdrs.org",
    "eurjther.com",
    "journals.iku.edu.tr",
    "www.jabsonline.org",
    "injectormedicaljournal.com",
    "www.medscidiscovery.com",
    "ijcmbs.com",
    "saglikokuryazarligidergisi.com",
    "natprobiotech.com",
    "www.derleme.gen.tr"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Remove the now unused rows (19-39) that previously held extra data
$ws.Range("A19:A39").ClearContents()
